## a69_f31_b UPP Pachuca — actualización trimestral (Q3 2021 -> Q4 2021)
## "actualizacion de febrero hay un archivo mal"
##
## Replaces the single Q3 "Informe financiero" row with the Q4 data, drops the
## two rows that were filed in error, clears the now-unused hyperlink columns
## on the kept row, and swaps the explanatory note for the one the finance
## office actually wants published.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# ------------------------------------------------------------------
# 1. Drop every hyperlink on the sheet (F8:F10, G8:G10). The three rows
#    get collapsed down to one, and that surviving row's link cells are
#    blanked out below, so none of the links stay relevant.
# ------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ------------------------------------------------------------------
# 2. The "Presupuesto Anual de Egresos Modificado" (row 9) and
#    "Adecuaciones Presupuestarias ..." (row 10) entries were the
#    mis-filed rows mentioned in the commit message — remove them
#    entirely, leaving only row 8.
# ------------------------------------------------------------------
$ws.Rows("9:10").Delete()

# ------------------------------------------------------------------
# 3. Update the remaining report row (row 8) to the 4th-trimester 2021
#    figures, matching the new reporting period 2021-10-01..2021-12-31,
#    validated/updated 2022-01-10.
# ------------------------------------------------------------------
$ws.Range("B8").Value = 44470   # 2021-10-01 fecha de inicio
$ws.Range("C8").Value = 44561   # 2021-12-31 fecha de termino

# Tipo de documento / denominacion / hipervinculos no longer apply to this
# row now that the two detail rows are gone — clear them out.
$ws.Range("D8:G8").ClearContents()
$ws.Range("D8:G8").HorizontalAlignment = -4131   # xlLeft

# "Área responsable" keeps its text, just re-left-align to match the rest
# of the row's cleared cells.
$ws.Range("H8").HorizontalAlignment = -4131      # xlLeft

$ws.Range("I8").Value = 44571   # fecha de validacion 2022-01-10
$ws.Range("J8").Value = 44571   # fecha de actualizacion 2022-01-10

# New note text explaining the late, definitive-figures publication.
$ws.Range("K8").Value = 'Derivado del Cierre financiero - presupuestal que se trabaja de manera  coordinada entre planeación y administración para la entrega y preparación de información ante las diversas dependencias fiscalizadoras , en apego al artículo 15 de la ley de fiscalización superior y rendición de cuenta de estado de hidalgo,  correlativamente con la fracción V del artículo 28 de la misma ley se establecen como fecha de entrega  los siete días hábiles siguientes al cierre del trimestre.  Así también dentro del convenio especifico para la  asignación de recursos con carácter de apoyo solidario firmado entre la federación y el estado en su cláusula sexta fracción "f" donde obliga la entrega de los estados financieros dentro de los primeros diez días hábiles a la coordinación de universidades tecnológicas y politécnicas. Motivo por el cual estaremos entregando la información con cifras definitivas, el 28 de enero del año en curso.'

# The note is long, so the row needs to grow to show it.
$ws.Rows(8).RowHeight = 208.5

# ------------------------------------------------------------------
# 4. The catalogue dropdown on column D only needs to stretch across the
#    (now much shorter) data range.
# ------------------------------------------------------------------
$ws.Range("D8:D201").Validation.Delete()
$ws.Range("D8:D113").Validation.Add(3, 1, 1, "Hidden_13")
$ws.Range("D8:D113").Validation.IgnoreBlank = $true
$ws.Range("D8:D113").Validation.InCellDropdown = $true
$ws.Range("D8:D113").Validation.ShowInput = $false
$ws.Range("D8:D113").Validation.ShowError = $true

# ------------------------------------------------------------------
# 5. Column widths: re-fit F/G now that the long hyperlink text is gone,
#    and trim K slightly to match the new note column.
# ------------------------------------------------------------------
$ws.Columns("F").ColumnWidth = 65.140625
$ws.Columns("G").ColumnWidth = 80.5703125
$ws.Columns("K").ColumnWidth = 71.42578125

# ------------------------------------------------------------------
# 6. Scroll/selection housekeeping to match where the editor left the
#    view after finishing the edit.
# ------------------------------------------------------------------
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 7
$ws.Range("B13").Select()
